# Testability Scenarios.docx - apply commit changes:
#  - remove the stray "_GoBack" bookmark after the author line
#  - prefix the three "Success percentage of the tests run" cells with "100% "
#    (split as "100% s" / "uccess percentage of the tests run", matching the
#    two-run shape produced by the original edit)
#  - tidy up a few places where the text was previously split across
#    multiple runs with no visible difference, by merging them back into a
#    single run ("Development, run time", "File System Listener module",
#    and the "Test that the system will work ... system." sentence)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark left over on the subtitle paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Prefix every "Success percentage of the tests run" with "100% ",
#    turning "Success ..." into "100% success ..." split over two runs:
#    "100% s" + "uccess percentage of the tests run".
# ---------------------------------------------------------------------
$needle = "Success percentage of the tests run"
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $mStart = $rng.Start
    $mEnd = $rng.End

    # Lower-case the leading "S" -> "s".
    $firstChar = $d.Range($mStart, $mStart + 1)
    $firstChar.Text = "s"

    # Insert "100% " right before it.
    $insertPoint = $d.Range($mStart, $mStart)
    $insertPoint.InsertBefore("100% ")

    # The inserted text + lower-cased "s" naturally merge into the run that
    # follows (same formatting), so force a split back into two runs by
    # toggling Bold on and back off across just the new "100% s" prefix.
    $splitRange = $d.Range($mStart, $insertPoint.End + 1)
    $splitRange.Bold = 1
    $splitRange.Bold = 0

    $searchStart = $insertPoint.End + ($mEnd - $mStart)
}

# ---------------------------------------------------------------------
# 3. Merge runs that are split with no actual formatting differences.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Development, run time", $true, $false, $false, $false, $false, $true, 1, $false, "Development, run time", 2) | Out-Null

$d.Content.Find.Execute("File System Listener module", $true, $false, $false, $false, $false, $true, 1, $false, "File System Listener module", 2) | Out-Null

$apos = [char]0x2019
$stimulus = "Test that the system will work with plugin files meant to work for another team" + $apos + "s system."
$d.Content.Find.Execute($stimulus, $true, $false, $false, $false, $false, $true, 1, $false, $stimulus, 2) | Out-Null
